# Update the fixed "Date" footer placeholder text from 2022/3/19 to 2022/8/21
# across the slide master, every slide layout, and the notes master.
$p = $ppt.ActivePresentation

$oldDate = "2022/3/19"
$newDate = "2022/8/21"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout (CustomLayout) attached to the master
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
